$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column widths for D, E, F (G keeps its existing bestFit width) ---
$ws.Columns.Item(4).ColumnWidth = 13
$ws.Columns.Item(5).ColumnWidth = 12.833333333333332
$ws.Columns.Item(6).ColumnWidth = 11.666666666666666

# --- Remove the now-unused explicit "applyNumberFormat" style from H2:O2 ---
$null = $ws.Range("H2:O2").ClearFormats()

# --- Add the new payroll row for Manu / STAFF ---
$ws.Range("A4").Value = "Manu"
$ws.Range("B4").Value = 1013
$ws.Range("C4").Value = "CSE"
$ws.Range("D4").Value = "STAFF"
$ws.Range("E4").Value = 44961
$ws.Range("F4").Value = 162063
$ws.Range("G4").Value = 5666009
$ws.Range("H4").Value = 67799
$ws.Range("I4").Value = 54444
$ws.Range("J4").Value = 676
$ws.Range("K4").Value = 7666
$ws.Range("L4").Value = 8999
$ws.Range("M4").Value = 788
$ws.Range("N4").Value = 89
$ws.Range("O4").Value = 988

# Give the new date cells (E4/F4) the same date number-format style already
# used by the other date column cells, reusing the existing style slot
# instead of minting a new one.
$null = $ws.Range("E2:F2").Copy()
$null = $ws.Range("E4:F4").PasteSpecial(-4122)

# --- Update view: zoom + selected cell ---
$ws.Application.ActiveWindow.Zoom = 66
$null = $ws.Range("Q9").Select()
